$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# --- New rows of expenses (20-22), already under construction ---

# Row 20: Geodeta cz.1
$ws.Range("A20").Value = "2013-07-xx"
$ws.Range("B20").Value = "Geodeta cz.1"
$ws.Range("B20").Font.Bold = $false
$ws.Range("C20").Value = 0
$ws.Range("C20").Font.Bold = $false
$ws.Range("D20").Value = 400
$ws.Range("E20").Formula = "=C20+D20"

# Row 21: Geodeta cz.2
$ws.Range("A21").Value = "2013-07-xx"
$ws.Range("B21").Value = "Geodeta cz.2"
$ws.Range("B21").Font.Bold = $false
$ws.Range("C21").Value = 0
$ws.Range("C21").Font.Bold = $false
$ws.Range("D21").Value = 650
$ws.Range("E21").Formula = "=C21+D21"

# Row 22: Piasek
$ws.Range("A22").Value = "2013-07-xx"
$ws.Range("B22").Value = "Piasek"
$ws.Range("B22").Font.Bold = $false
$ws.Range("C22").Value = 2000
$ws.Range("C22").Font.Bold = $false
$ws.Range("D22").Value = 5800
$ws.Range("E22").Formula = "=C22+D22"

# --- Update the view: scroll back to top and select the newly entered cell ---
$ws.Activate() | Out-Null
$ws.Range("B22").Select() | Out-Null

Write-Output "applied"
